$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Apply replacements in document order so that a newly-introduced value
# (e.g. "88÷7=" created from "64÷2=") is never re-matched by a later rule
# whose "old" text happens to equal that new value.
Replace-Text "76÷4=" "73÷3="
Replace-Text "56÷5=" "21÷7="
Replace-Text "77÷8=" "76÷2="
Replace-Text "88÷7=" "85÷4="
Replace-Text "30÷9=" "26÷5="
Replace-Text "95÷3=" "31÷3="
Replace-Text "80÷6=" "88÷5="
Replace-Text "58÷9=" "23÷4="
Replace-Text "21÷5=" "18÷7="
Replace-Text "12÷7=" "44÷6="
Replace-Text "66÷7=" "25÷6="
Replace-Text "85÷7=" "46÷8="
Replace-Text "22÷9=" "38÷9="
Replace-Text "55÷3=" "79÷3="
Replace-Text "12÷8=" "40÷2="
Replace-Text "87÷4=" "64÷8="
Replace-Text "78÷8=" "18÷4="
Replace-Text "33÷6=" "97÷8="
Replace-Text "37÷7=" "75÷7="
Replace-Text "64÷2=" "88÷7="
Replace-Text "60÷9=" "59÷5="
Replace-Text "78÷3=" "13÷2="
Replace-Text "17÷6=" "98÷5="
Replace-Text "46÷2=" "43÷7="
Replace-Text "77÷5=" "87÷7="

Write-Output "Replacements applied"
